$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3167.6667
$ws.Range("I74").Value = 3001.5
$ws.Range("K74").Value = 3001.5
$ws.Range("M74").Value = -2065.5
# Row 77
$ws.Range("H77").Value = 3167.6667
$ws.Range("I77").Value = 3001.5
$ws.Range("K77").Value = 15007.5
$ws.Range("M77").Value = -10327.5
# Row 112
$ws.Range("H112").Value = 923.4
$ws.Range("J112").Value = 923.4
$ws.Range("L112").Value = 2770.2
$ws.Range("N112").Value = -4986.2
# Row 138
$ws.Range("H138").Value = 5149.7256
$ws.Range("J138").Value = 6182.65
$ws.Range("L138").Value = 18547.95
$ws.Range("N138").Value = -28827.95

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1346.1428
$ws.Range("I2").Value = 1345.8334
$ws.Range("K2").Value = 1345.8334
$ws.Range("M2").Value = -1232.8334
# Row 32
$ws.Range("H32").Value = 3650
$ws.Range("I32").Value = 3630.5557
$ws.Range("K32").Value = 3630.5557
$ws.Range("M32").Value = -3343.5557
# Row 43
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50626
# Row 45
$ws.Range("H45").Value = 2258
$ws.Range("I45").Value = 2298
$ws.Range("J45").Value = 2164.6667
$ws.Range("K45").Value = 2298
$ws.Range("L45").Value = 2164.6667
$ws.Range("M45").Value = -1921
$ws.Range("N45").Value = -2918.6667
# Row 74
$ws.Range("H74").Value = 2733.5217
$ws.Range("I74").Value = 633.6429000000001
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 633.6429000000001
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = 240.3570999999999
$ws.Range("N74").Value = -7748
# Row 77
$ws.Range("H77").Value = 2733.5217
$ws.Range("I77").Value = 633.6429000000001
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 3168.2145
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = 1199.7855
$ws.Range("N77").Value = -38736
# Row 97
$ws.Range("H97").Value = 806
$ws.Range("I97").Value = 661
$ws.Range("J97").Value = 1096
$ws.Range("K97").Value = 661
$ws.Range("L97").Value = 1096
$ws.Range("M97").Value = -165
$ws.Range("N97").Value = -2088
# Row 116
$ws.Range("H116").Value = 1346.1428
$ws.Range("I116").Value = 1345.8334
$ws.Range("K116").Value = 1345.8334
$ws.Range("M116").Value = 948.1666
# Row 139
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1346.1428
$ws.Range("I3").Value = 1345.8334
$ws.Range("K3").Value = 1345.8334
$ws.Range("M3").Value = -1231.8334
# Row 80
$ws.Range("H80").Value = 529.36365
$ws.Range("I80").Value = 1542.6666
$ws.Range("J80").Value = 149.375
$ws.Range("K80").Value = 1542.6666
$ws.Range("L80").Value = 149.375
$ws.Range("M80").Value = -544.6666
$ws.Range("N80").Value = -2145.375
# Row 83
$ws.Range("H83").Value = 529.36365
$ws.Range("I83").Value = 1542.6666
$ws.Range("J83").Value = 149.375
$ws.Range("K83").Value = 7713.333000000001
$ws.Range("L83").Value = 746.875
$ws.Range("M83").Value = -2721.333000000001
$ws.Range("N83").Value = -10730.875
# Row 86
$ws.Range("H86").Value = 3117.8
$ws.Range("I86").Value = 3196.3333
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3196.3333
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -2073.3333
$ws.Range("N86").Value = -5246
# Row 89
$ws.Range("H89").Value = 3117.8
$ws.Range("I89").Value = 3196.3333
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15981.6665
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -10365.6665
$ws.Range("N89").Value = -26232
# Row 94
$ws.Range("H94").Value = 1010.75
$ws.Range("I94").Value = 817.4
$ws.Range("J94").Value = 1333
$ws.Range("K94").Value = 817.4
$ws.Range("L94").Value = 1333
$ws.Range("M94").Value = -366.4
$ws.Range("N94").Value = -2235
# Row 105
$ws.Range("H105").Value = 2856.5
$ws.Range("I105").Value = 2856.5
$ws.Range("K105").Value = 2856.5
$ws.Range("M105").Value = -1109.5
# Row 107
$ws.Range("H107").Value = 5927.6665
$ws.Range("I107").Value = 5941.5
$ws.Range("J107").Value = 5900
$ws.Range("K107").Value = 5941.5
$ws.Range("L107").Value = 5900
$ws.Range("M107").Value = -4021.5
$ws.Range("N107").Value = -9740
# Row 138
$ws.Range("H138").Value = 124499
$ws.Range("J138").Value = 124499
$ws.Range("L138").Value = 124499
$ws.Range("N138").Value = -134779

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 8163.3335
$ws.Range("J41").Value = 7993.3335
$ws.Range("L41").Value = 7993.3335
$ws.Range("N41").Value = -8849.333500000001
# Row 122
$ws.Range("H122").Value = 1595.6
$ws.Range("I122").Value = 1595.6
$ws.Range("K122").Value = 4786.799999999999
$ws.Range("M122").Value = -2336.799999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1188.5
$ws.Range("I97").Value = 765
$ws.Range("K97").Value = 765
$ws.Range("M97").Value = -269
# Row 132
$ws.Range("H132").Value = 5998.3335
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1832.5
$ws.Range("I16").Value = 1999
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1999
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1829
$ws.Range("N16").Value = -1340
# Row 40
$ws.Range("H40").Value = 4260.3
$ws.Range("I40").Value = 4370.143
$ws.Range("J40").Value = 4004
$ws.Range("K40").Value = 4370.143
$ws.Range("L40").Value = 4004
$ws.Range("M40").Value = -4234.143
$ws.Range("N40").Value = -4276
# Row 68
$ws.Range("H68").Value = 2984.6365
$ws.Range("I68").Value = 2984.6365
$ws.Range("K68").Value = 2984.6365
$ws.Range("M68").Value = -2235.6365
# Row 71
$ws.Range("H71").Value = 2984.6365
$ws.Range("I71").Value = 2984.6365
$ws.Range("K71").Value = 14923.1825
$ws.Range("M71").Value = -11179.1825
# Row 93
$ws.Range("H93").Value = 499
$ws.Range("I93").Value = 499
$ws.Range("K93").Value = 499
$ws.Range("M93").Value = 749
# Row 132
$ws.Range("H132").Value = 4999.3335
# Row 136
$ws.Range("H136").Value = 21660.715
$ws.Range("I136").Value = 10325.4
$ws.Range("K136").Value = 30976.2
$ws.Range("M136").Value = -28426.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 9539.25
$ws.Range("I81").Value = 6558.375
$ws.Range("J81").Value = 15501
$ws.Range("K81").Value = 13116.75
$ws.Range("L81").Value = 31002
$ws.Range("M81").Value = -12055.75
$ws.Range("N81").Value = -33124
# Row 84
$ws.Range("H84").Value = 9539.25
$ws.Range("I84").Value = 6558.375
$ws.Range("J84").Value = 15501
$ws.Range("K84").Value = 65583.75
$ws.Range("L84").Value = 155010
$ws.Range("M84").Value = -60279.75
$ws.Range("N84").Value = -165618
# Row 113
$ws.Range("H113").Value = 1013.5714
$ws.Range("I113").Value = 919.2
$ws.Range("J113").Value = 1249.5
$ws.Range("K113").Value = 2757.6
$ws.Range("L113").Value = 3748.5
$ws.Range("M113").Value = -587.6000000000004
$ws.Range("N113").Value = -8088.5
# Row 132
$ws.Range("H132").Value = 3409.162
$ws.Range("I132").Value = 2864.625
$ws.Range("J132").Value = 4414.4614
$ws.Range("K132").Value = 8593.875
$ws.Range("L132").Value = 13243.3842
$ws.Range("M132").Value = -6063.875
$ws.Range("N132").Value = -18303.3842
